# advertising-packets.xlsx update:
#  - add a second (empty) worksheet after Sheet1
#  - rework the "packet type + setting" header block (G1:J1 / G4:J4)
#  - shift the old bit-offset column from E to F and renumber it
#  - replace the lsm9 calibration block (old H7:I10 + L10 formulas) with a
#    name/value "setting" table in H:J that runs all the way to row 24
#  - extend the sheet down to row 27 with the new bit-offset numbering
#
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- add Sheet2 (placed after Sheet1, like the real edit) ----------------
[void]$wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)

# ---- row 1 / row 2 header cells ------------------------------------------
$ws.Range("C1").Value = 0
$ws.Range("G1").Value = 0
$ws.Range("H1").Value = 1
$ws.Range("I1").Value = 2
$ws.Range("J1").Value = 3

# ---- drop the old "bits" column E (rows 2-19) -----------------------------
$ws.Range("E2:E19").Clear()

# ---- drop the old lsm9-calibration formulas/values ------------------------
$ws.Range("I6").Clear()
$ws.Range("H7:I9").Clear()
$ws.Range("I10").Clear()
$ws.Range("L10").Clear()

# ---- F column: renumbered bit offsets (row2 -> 0, row3 -> 1, ... row27 -> 25)
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 2
$ws.Range("F5").Value = 3
$ws.Range("F6").Value = 4
$ws.Range("F7").Value = 5
$ws.Range("F8").Value = 6
$ws.Range("F9").Value = 7
$ws.Range("F10").Value = 8
$ws.Range("F11").Value = 9
$ws.Range("F12").Value = 10
$ws.Range("F13").Value = 11
$ws.Range("F14").Value = 12
$ws.Range("F15").Value = 13
$ws.Range("F16").Value = 14
$ws.Range("F17").Value = 15
$ws.Range("F18").Value = 16
$ws.Range("F19").Value = 17
$ws.Range("F20").Value = 18
$ws.Range("F21").Value = 19
$ws.Range("F22").Value = 20
$ws.Range("F23").Value = 21
$ws.Range("F24").Value = 22
$ws.Range("F25").Value = 23
$ws.Range("F26").Value = 24
$ws.Range("F27").Value = 25

# ---- new "packet type + setting" header block (row 4) --------------------
$ws.Range("G4").Value = "packet type + setting"
$ws.Range("H4").Value = "packet type"
$ws.Range("I4").Value = "packet type"
$ws.Range("J4").Value = "packet type"

# ---- J column: setting "name" label, top to bottom ------------------------
$ws.Range("J5:J18").Value = "name"
$ws.Range("J19:J22").Value = "wheel"
$ws.Range("J23:J24").Value = "zero pos"

# ---- H / I columns: settings name/value pairs, top to bottom -------------
$ws.Range("H5").Value = "tb0"
$ws.Range("I5").Value = "tb0"

$ws.Range("H6").Value = "tb1"
$ws.Range("I6").Value = "tb1"

$ws.Range("H7").Value = "tb2"
$ws.Range("I7").Value = "tb2"

$ws.Range("H8:H11").Value = "humidity"
$ws.Range("I8:I11").Value = "altitude"

$ws.Range("H12:H15").Value = "air pressure"
$ws.Range("I12:I15").Value = "air density"

$ws.Range("H16:H19").Value = "temp"
$ws.Range("I16:I19").Value = "dew point"

# ---- move the active selection the same place the author left it ----------
[void]$ws.Range("F18").Select()

# ---- column G gets a bestFit-like width like column B already has --------
$ws.Columns.Item(7).ColumnWidth = 18.1640625
